# Atualizado por script em 02-12-2023 14:46
#
# This script reorders several match rows that share the same matchday
# (same "data_partida") back into the order the source site now reports
# them in, and appends 3 newly scraped matches (rows 161-163) for the
# 2023-12-02 matchday.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Re-shuffle columns F:V (home .. url_partida) among rows that share
#        the same Indice/pais/torneio/temporada/data_partida (columns A:E,
#        left untouched). $mapping[$r] = $src means "row $r should end up
#        holding the F:V content that currently lives in row $src".
$mapping = @{
    12 = 13;  13 = 12;
    24 = 25;  25 = 24;
    29 = 31;  30 = 32;  31 = 30;  32 = 29;
    37 = 39;  38 = 40;  39 = 38;  40 = 37;
    73 = 74;  74 = 75;  75 = 76;  76 = 77;  77 = 73;
    82 = 83;  83 = 82;
    130 = 132; 132 = 130;
    144 = 145; 145 = 144;
    146 = 149; 147 = 150; 148 = 147; 149 = 146; 150 = 148;
}

# Snapshot the current F:V content of every affected row first so the
# subsequent writes can't clobber a value that still needs to be read.
$snapshot = @{}
foreach ($r in $mapping.Keys) {
    $snapshot[$r] = $ws.Range("F${r}:V${r}").Value()
}

foreach ($r in $mapping.Keys) {
    $src = $mapping[$r]
    $ws.Range("F${r}:V${r}").Value = $snapshot[$src]
}

# --- 2) Append the 3 newly scraped rows (161-163), matching the existing
#        formatting (bold/bordered index column, date-formatted E column).
$ws.Range("A160:V160").Copy()
$ws.Range("A161:V163").PasteSpecial(-4122)

$vals161 = @(160, "bulgaria", "vtora-liga", "2023-2024", 45262.54166666666, `
    "Dunav Ruse", 0, "Spartak Pleven", 1, `
    1.49, "02/12/2023 02:13", 1.44, "02/12/2023 12:59", `
    3.76, "02/12/2023 02:13", 4.15, "02/12/2023 12:59", `
    5.67, "02/12/2023 02:13", 5.88, "02/12/2023 12:59", `
    "https://www.betexplorer.com/football/bulgaria/vtora-liga/dunav-ruse-spartak-pleven/4j2Hvnot/")
$arr161 = New-Object 'object[,]' 1,22
for ($i = 0; $i -lt 22; $i++) { $arr161[0,$i] = $vals161[$i] }
$ws.Range("A161:V161").Value = $arr161

$vals162 = @(161, "bulgaria", "vtora-liga", "2023-2024", 45262.54166666666, `
    "Belasitsa", 0, "Yantra Gabrovo", 0, `
    3.01, "02/12/2023 02:13", 2.29, "02/12/2023 12:57", `
    2.78, "02/12/2023 02:13", 2.81, "02/12/2023 12:57", `
    2.37, "02/12/2023 02:13", 3.2, "02/12/2023 12:57", `
    "https://www.betexplorer.com/football/bulgaria/vtora-liga/belasitsa-petrich-yantra-gabrovo/j74TypVb/")
$arr162 = New-Object 'object[,]' 1,22
for ($i = 0; $i -lt 22; $i++) { $arr162[0,$i] = $vals162[$i] }
$ws.Range("A162:V162").Value = $arr162

$vals163 = @(162, "bulgaria", "vtora-liga", "2023-2024", 45262.54166666666, `
    "Strumska Slava", 2, "Chernomorets Balchik", 2, `
    1.48, "02/12/2023 02:13", 1.35, "02/12/2023 12:00", `
    3.74, "02/12/2023 02:13", 4.08, "02/12/2023 12:55", `
    5.98, "02/12/2023 02:13", 8.369999999999999, "02/12/2023 12:55", `
    "https://www.betexplorer.com/football/bulgaria/vtora-liga/strumska-slava-chernomorets-balchik/hCx3LOpI/")
$arr163 = New-Object 'object[,]' 1,22
for ($i = 0; $i -lt 22; $i++) { $arr163[0,$i] = $vals163[$i] }
$ws.Range("A163:V163").Value = $arr163
